# Natmi following Dr Hou advice:
# recompute the Inhba-Acvr2b ligand-receptor table for all 4x4
# sending/target cluster combinations (ECs, FAPs, M2, sCs) with the
# revised expression statistics, extending the table from 12 to 16
# data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @('ECs', 'Inhba', 'Acvr2b', 'ECs', 3, 1, 2.288150666666667, 6.864452, 0.3964219041944151, 0.3964219041944151, 2, 0.6666666666666666, 1.047294333333333, 3.141883, 0.2943898468263134, 0.2943898468263134, 2.396367227012889, 21.567305043116, 0.1167025836543893, 0.1167025836543893),
  @('ECs', 'Inhba', 'Acvr2b', 'FAPs', 3, 1, 2.288150666666667, 6.864452, 0.3964219041944151, 0.3964219041944151, 3, 1, 0.5210533333333333, 1.56316, 0.1464658082318851, 0.1464658082318851, 1.192248532035555, 10.73023678832, 0.05806225459865794, 0.05806225459865794),
  @('ECs', 'Inhba', 'Acvr2b', 'M2', 3, 1, 2.288150666666667, 6.864452, 0.3964219041944151, 0.3964219041944151, 3, 1, 0.3665986666666667, 1.099796, 0.1030492784041265, 0.1030492784041265, 0.8388329835324445, 7.549496851792, 0.04085099117082424, 0.04085099117082423),
  @('ECs', 'Inhba', 'Acvr2b', 'sCs', 3, 1, 2.288150666666667, 6.864452, 0.3964219041944151, 0.3964219041944151, 3, 1, 1.622562, 4.867686, 0.456095066537675, 0.456095066537675, 3.712666322008, 33.413996898072, 0.1808060747705436, 0.1808060747705436),
  @('FAPs', 'Inhba', 'Acvr2b', 'ECs', 3, 1, 2.588894, 7.766681999999999, 0.4485256605643812, 0.4485256605643813, 2, 0.6666666666666666, 1.047294333333333, 3.141883, 0.2943898468263134, 0.2943898468263134, 2.711334015800666, 24.402006142206, 0.1320414005112192, 0.1320414005112192),
  @('FAPs', 'Inhba', 'Acvr2b', 'FAPs', 3, 1, 2.588894, 7.766681999999999, 0.4485256605643812, 0.4485256605643813, 3, 1, 0.5210533333333333, 1.56316, 0.1464658082318851, 0.1464658082318851, 1.348951848346666, 12.14056663512, 0.06569367338730225, 0.06569367338730227),
  @('FAPs', 'Inhba', 'Acvr2b', 'M2', 3, 1, 2.588894, 7.766681999999999, 0.4485256605643812, 0.4485256605643813, 3, 1, 0.3665986666666667, 1.099796, 0.1030492784041265, 0.1030492784041265, 0.9490850885413333, 8.541765796871999, 0.04622024566689366, 0.04622024566689366),
  @('FAPs', 'Inhba', 'Acvr2b', 'sCs', 3, 1, 2.588894, 7.766681999999999, 0.4485256605643812, 0.4485256605643813, 3, 1, 1.622562, 4.867686, 0.456095066537675, 0.456095066537675, 4.200641026427999, 37.805769237852, 0.2045703409989661, 0.2045703409989661),
  @('M2', 'Inhba', 'Acvr2b', 'ECs', 3, 1, 0.3123523333333333, 0.9370569999999999, 0.05411501461132016, 0.05411501461132018, 2, 0.6666666666666666, 1.047294333333333, 3.141883, 0.2943898468263134, 0.2943898468263134, 0.3271248287034444, 2.944123458331, 0.01593091086243025, 0.01593091086243026),
  @('M2', 'Inhba', 'Acvr2b', 'FAPs', 3, 1, 0.3123523333333333, 0.9370569999999999, 0.05411501461132016, 0.05411501461132018, 3, 1, 0.5210533333333333, 1.56316, 0.1464658082318851, 0.1464658082318851, 0.1627522244577777, 1.46477002012, 0.007925999352527281, 0.007925999352527283),
  @('M2', 'Inhba', 'Acvr2b', 'M2', 3, 1, 0.3123523333333333, 0.9370569999999999, 0.05411501461132016, 0.05411501461132018, 3, 1, 0.3665986666666667, 1.099796, 0.1030492784041265, 0.1030492784041265, 0.1145079489302222, 1.030571540372, 0.005576513206525305, 0.005576513206525305),
  @('M2', 'Inhba', 'Acvr2b', 'sCs', 3, 1, 0.3123523333333333, 0.9370569999999999, 0.05411501461132016, 0.05411501461132018, 3, 1, 1.622562, 4.867686, 0.456095066537675, 0.456095066537675, 0.5068110266779999, 4.561299240102, 0.02468159118983733, 0.02468159118983733),
  @('sCs', 'Inhba', 'Acvr2b', 'ECs', 3, 1, 0.5826116666666666, 1.747835, 0.1009374206298835, 0.1009374206298836, 2, 0.6666666666666666, 1.047294333333333, 3.141883, 0.2943898468263134, 0.2943898468263134, 0.6101658970338889, 5.491493073305, 0.02971495179827458, 0.02971495179827459),
  @('sCs', 'Inhba', 'Acvr2b', 'FAPs', 3, 1, 0.5826116666666666, 1.747835, 0.1009374206298835, 0.1009374206298836, 3, 1, 0.5210533333333333, 1.56316, 0.1464658082318851, 0.1464658082318851, 0.3035717509555555, 2.732145758599999, 0.01478388089339765, 0.01478388089339765),
  @('sCs', 'Inhba', 'Acvr2b', 'M2', 3, 1, 0.5826116666666666, 1.747835, 0.1009374206298835, 0.1009374206298836, 3, 1, 0.3665986666666667, 1.099796, 0.1030492784041265, 0.1030492784041265, 0.2135846601844444, 1.92226194166, 0.01040152835988329, 0.01040152835988329),
  @('sCs', 'Inhba', 'Acvr2b', 'sCs', 3, 1, 0.5826116666666666, 1.747835, 0.1009374206298835, 0.1009374206298836, 3, 1, 1.622562, 4.867686, 0.456095066537675, 0.456095066537675, 0.94532355109, 8.507911959809999, 0.04603705957832803, 0.04603705957832804)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}